$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (D) and Volume 1h (E) columns for changed rows ---
$ws.Range("D2").Value = "29.431.92"
$ws.Range("E2").Value = "  +0.43%  "

$ws.Range("D3").Value = "1.848.87"
$ws.Range("E3").Value = "  +0.43%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "240.79"
$ws.Range("E5").Value = "  +0.82%  "

$ws.Range("D6").Value = "0.6293"
$ws.Range("E6").Value = "  +0.08%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").Value = "0.07696"
$ws.Range("E8").Value = "  +2.33%  "

$ws.Range("E9").Value = "  -0.45%  "

$ws.Range("D10").Value = "24.68"
$ws.Range("E10").Value = "  +0.96%  "

$ws.Range("E11").Value = "  +0.71%  "

$ws.Range("D12").Value = "1.886.08"
$ws.Range("E12").Value = "  +2.16%  "

$ws.Range("D13").Value = "5.031"
$ws.Range("E13").Value = "  +1.12%  "

$ws.Range("D14").Value = "0.00001080"
$ws.Range("E14").Value = "  +4.17%  "

$ws.Range("D15").Value = "0.6787"
$ws.Range("E15").Value = "  +0.25%  "

$ws.Range("D16").Value = "83.61"
$ws.Range("E16").Value = "  +0.73%  "

$ws.Range("D17").Value = "2.147.88"
$ws.Range("E17").Value = "  +2.05%  "

$ws.Range("D18").Value = "6.194"
$ws.Range("E18").Value = "  +0.94%  "

$ws.Range("D19").Value = "29.470.62"
$ws.Range("E19").Value = "  +0.42%  "

$ws.Range("D20").Value = "228.02"
$ws.Range("E20").Value = "  -0.16%  "

$ws.Range("D21").Value = "12.44"
$ws.Range("E21").Value = "  +0.36%  "

$ws.Range("E22").Value = "  +0.01%  "

$ws.Range("D23").Value = "7.431"
$ws.Range("E23").Value = "  -0.04%  "

$ws.Range("E24").Value = "  +0.01%  "

$ws.Range("D25").Value = "157.61"
$ws.Range("E25").Value = "  +0.77%  "

$ws.Range("D26").Value = "0.1380"
$ws.Range("E26").Value = "  -0.75%  "

$ws.Range("E27").Value = "  +0.81%  "

$ws.Range("E28").Value = "  +0.55%  "

$ws.Range("D29").Value = "1.345"
$ws.Range("E29").Value = "  +5.84%  "

$ws.Range("D30").Value = "1.469"
$ws.Range("E30").Value = "  +0.69%  "

$ws.Range("E31").Value = "  +0.60%  "

$ws.Range("D32").Value = "4.125"

$ws.Range("D33").Value = "4.035"

$ws.Range("E34").Value = "  +1.09%  "

$ws.Range("D35").Value = "1.162"
$ws.Range("E35").Value = "  +0.75%  "

$ws.Range("D36").Value = "0.7018"
$ws.Range("E36").Value = "  -0.98%  "

$ws.Range("E37").Value = "  -0.34%  "

$ws.Range("E38").Value = "  +0.50%  "

$ws.Range("E39").Value = "  -0.95%  "

$ws.Range("D40").Value = "1.220.41"
$ws.Range("E40").Value = "  -1.62%  "

$ws.Range("D41").Value = "6.541"
$ws.Range("E41").Value = "  +5.05%  "

$ws.Range("D42").Value = "0.9055"
$ws.Range("E42").Value = "  +0.50%  "

$ws.Range("E43").Value = "  +0.11%  "

$ws.Range("D44").Value = "101.80"
$ws.Range("E44").Value = "  +0.08%  "

$ws.Range("D45").Value = "66.26"

# --- Rows 46 and 47 swapped content (Aptos <-> BabyDogeCoin) ---
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.00000000120"
$ws.Range("E46").Value = "  +1.60%  "

$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").Value = "7.137"
$ws.Range("E47").Value = "  +0.60%  "


$ws.Range("D48").Value = "0.4029"
$ws.Range("E48").Value = "  +1.01%  "

$ws.Range("D49").Value = "9.054"
$ws.Range("E49").Value = "  +1.20%  "

$ws.Range("E50").Value = "  +0.85%  "

$ws.Range("D51").Value = "0.1145"
$ws.Range("E51").Value = "  +2.22%  "
